# Support const char * argument type and size_t return type.
# Adds a new "MultiByteStrLen" worked example (2 rows) between the
# "ReverseString" example and the "Array Example" / Trace example,
# and moves the selection to D13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a single row at row 12 - this pushes the old row 14 ("Array
# Example" label) and everything below it down by one row (14->15,
# 15->16, 16->17, 17->18, 18->19), and the new row 12 inherits the
# formatting of row 11 (the row immediately above the insertion point).
$ws.Rows.Item(12).Insert()

# ---- Row 12: MultiByteStrLen("hello") -> 5 ----
$ws.Range("B12").Value = "MultiByteStrLen"
$ws.Range("C12").Value = "hello"
$ws.Range("F12").Formula = "=_xll.MultiByteStrLen(C12)"
$ws.Range("G12").Value = 5
$ws.Range("H12").Formula = "=F12=G12"

# ---- Row 13: MultiByteStrLen(REPT("x",256)) -> error ----
$ws.Range("B13").Value = "MultiByteStrLen"
$ws.Range("C13").Value = "x"
$ws.Range("D13").Value = 256
$ws.Range("F13").Formula = "=_xll.MultiByteStrLen(REPT(C13,D13))"
$ws.Range("G13").Formula = "=LEN(C13)*D13"
$ws.Range("H13").Formula = "=F13=G13"

# Row 13 was not created by the row-insert shift, so it does not
# automatically inherit row 12's look; copy the formatting over.
$ws.Range("B12:H12").Copy()
$ws.Range("B13:H13").PasteSpecial(-4122)

# Match the workbook's recorded selection after the edit.
$ws.Range("D13").Select()
